$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: capture current N1 style (bordered empty-header style) onto P1 first,
#     before we touch N1's own style. ---
$ws.Range("N1").Copy()
$ws.Range("P1").PasteSpecial(-4122)

# --- Step 2: upgrade N1 (currently style s=3) and prep O1 to the data-header
#     style s=4 by copying format from a cell that already has it (M1). ---
$ws.Range("M1").Copy()
$ws.Range("N1:O1").PasteSpecial(-4122)

# --- Step 3: write the new header text, in final column order. ---
$ws.Range("A1").Value = "Скважина"
$ws.Range("B1").Value = "Дата"
$ws.Range("C1").Value = "Нефть факт"
$ws.Range("D1").Value = "Жидкость факт"
$ws.Range("E1").Value = "Отработанные часы"
$ws.Range("F1").Value = "Часы в простое"
$ws.Range("G1").Value = "Состояние"
$ws.Range("H1").Value = "Доли ПРС"
$ws.Range("I1").Value = "Причина потерь"
$ws.Range("J1").Value = "Нефть прогноз"
$ws.Range("K1").Value = "Потери нефти"
$ws.Range("L1").Value = "Жидкость прогноз"
$ws.Range("M1").Value = "Потери жидкости"
$ws.Range("N1").Value = "Тех потери нефти"
$ws.Range("O1").Value = "Тех потери жидкости"

# --- Step 4: row height back to "auto" (drops ht/customHeight like the diff). ---
$ws.Rows.Item(1).AutoFit()

# --- Step 5: column widths (best effort within the host's width quantisation). ---
$ws.Columns.Item(1).ColumnWidth = 10.333333333333334
$ws.Columns.Item(2).ColumnWidth = 10.333333333333334
$ws.Columns.Item(3).ColumnWidth = 10.333333333333334
$ws.Columns.Item(4).ColumnWidth = 10.333333333333334
$ws.Columns.Item(5).ColumnWidth = 20.666666666666668
$ws.Columns.Item(6).ColumnWidth = 12.666666666666666
$ws.Columns.Item(7).ColumnWidth = 10.333333333333334
$ws.Columns.Item(8).ColumnWidth = 10.333333333333334
$ws.Columns.Item(9).ColumnWidth = 15.5
$ws.Columns.Item(10).ColumnWidth = 10.333333333333334
$ws.Columns.Item(11).ColumnWidth = 11.5
$ws.Columns.Item(12).ColumnWidth = 18.166666666666668
$ws.Columns.Item(13).ColumnWidth = 15.666666666666666
$ws.Columns.Item(14).ColumnWidth = 15.0

# --- Step 6: selection. ---
$ws.Range("J4").Select()
